$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Add new row 63 with the new time entry
$ws.Cells.Item(63, 1).Value = "Federico Speroni"

$ws.Cells.Item(63, 2).Value = 42893
$ws.Range("B62").Copy()
$ws.Range("B63").PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item(63, 3).Value = 3
$ws.Cells.Item(63, 4).Value = "Sprint 3 - FrontEnd"
$ws.Cells.Item(63, 5).Value = "Administrador - Arreglos en funcionalidades. Algunas Pruebas IU cliente y administrador"

# Update selection to match the new last row
$ws.Range("E63").Select()

$wb.Save()
